$d = $word.ActiveDocument

# Update the date heading
$dateFound = $d.Content.Find.Execute("2025-05-30 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-05-31 Saturday", 2)
if (-not $dateFound) {
    Write-Host "WARNING: date heading text not found"
}

# Update the multiplication problems in the table, cell by cell so that
# replacement values that coincide with other cells' original text never
# collide with each other.
$t = $d.Tables.Item(1)

$cellUpdates = @(
    @{ Row = 1;  Col = 1; Old = "21×40="; New = "54×27=" },
    @{ Row = 1;  Col = 2; Old = "91×64="; New = "25×91=" },
    @{ Row = 1;  Col = 3; Old = "84×58="; New = "57×19=" },
    @{ Row = 1;  Col = 4; Old = "84×75="; New = "84×37=" },
    @{ Row = 1;  Col = 5; Old = "19×52="; New = "68×64=" },

    @{ Row = 5;  Col = 1; Old = "53×49="; New = "90×28=" },
    @{ Row = 5;  Col = 2; Old = "82×90="; New = "33×11=" },
    @{ Row = 5;  Col = 3; Old = "90×44="; New = "84×83=" },
    @{ Row = 5;  Col = 4; Old = "76×91="; New = "98×64=" },
    @{ Row = 5;  Col = 5; Old = "49×21="; New = "77×73=" },

    @{ Row = 10; Col = 1; Old = "88×23="; New = "72×29=" },
    @{ Row = 10; Col = 2; Old = "85×96="; New = "17×91=" },
    @{ Row = 10; Col = 3; Old = "93×27="; New = "37×12=" },
    @{ Row = 10; Col = 4; Old = "75×62="; New = "63×58=" },
    @{ Row = 10; Col = 5; Old = "88×64="; New = "51×69=" },

    @{ Row = 15; Col = 1; Old = "58×84="; New = "29×96=" },
    @{ Row = 15; Col = 2; Old = "95×60="; New = "21×40=" },
    @{ Row = 15; Col = 3; Old = "67×46="; New = "79×20=" },
    @{ Row = 15; Col = 4; Old = "79×72="; New = "16×61=" },
    @{ Row = 15; Col = 5; Old = "50×61="; New = "49×63=" },

    @{ Row = 20; Col = 1; Old = "71×72="; New = "83×63=" },
    @{ Row = 20; Col = 2; Old = "75×87="; New = "26×28=" },
    @{ Row = 20; Col = 3; Old = "17×13="; New = "46×90=" },
    @{ Row = 20; Col = 4; Old = "75×57="; New = "21×99=" },
    @{ Row = 20; Col = 5; Old = "67×81="; New = "31×99=" }
)

foreach ($u in $cellUpdates) {
    $cellRange = $t.Cell($u.Row, $u.Col).Range
    $found = $cellRange.Find.Execute($u.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $u.New, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found in cell" $u.Row $u.Col ":" $u.Old
    }
}
